$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 10613.61331316024
$ws.Range("D3").Value = 697.1794862480978

$ws.Range("B4").Value = 4446.758018384022
$ws.Range("D4").Value = 538.3281476284172

$ws.Range("B5").Value = 3238.336364383562
$ws.Range("D5").Value = 27.94551232876716

$ws.Range("B6").Value = 7017.177026712333
$ws.Range("D6").Value = 158.3569445205474

$ws.Range("B7").Value = 9203.445441095897
$ws.Range("D7").Value = 657.3151780821913

$ws.Range("B8").Value = 14316.09068561651
$ws.Range("D8").Value = 1280

$ws.Range("B9").Value = 20122.64999452061
$ws.Range("D9").Value = 1284.657580821918

$ws.Range("F10").Value = 15283603.17095072

$ws.Range("G11").Value = 0.7210308521362604

$ws.Range("F12").Value = 1003938.46019726
$ws.Range("G12").Value = 0.06568728911422068

$ws.Range("G13").Value = 0.2132818587495189
